# Update the "Fuel" column (I) values for rows 2-21 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$fuelValues = @{
    2  = 25
    3  = 26
    4  = 31
    5  = 42
    6  = 11
    7  = 23
    8  = 56
    9  = 74
    10 = 35
    11 = 34
    12 = 45
    13 = 12
    14 = 65
    15 = 43
    16 = 24
    17 = 53
    18 = 27
    19 = 75
    20 = 78
    21 = 77
}

foreach ($row in $fuelValues.Keys) {
    $ws.Range("I$row").Value = $fuelValues[$row]
}

# Update the column widths to reflect the resaved ("best fit" recalculated) layout.
$ws.Columns.Item(1).ColumnWidth = 13.998697916666666
$ws.Columns.Item(2).ColumnWidth = 15.830729166666666
$ws.Columns.Item(3).ColumnWidth = 39.830729166666664
$ws.Columns.Item(4).ColumnWidth = 36.498697916666664
$ws.Columns.Item(5).ColumnWidth = 46.498697916666664

# Update the selected cell on the sheet (matches the new <selection> entry).
$ws.Range("I22").Select()
